# Remove footnote markers like " [1]" / " [5, 6]" from vaccine/category names,
# and collapse any embedded line breaks in cell text into a single space.
# (Also fixes a data-entry inconsistency on the "Adult Influenza Vaccine"
#  sheet where some rows used "Influenza [5, 6]" and others "Influenza [5]"
#  for the same column - after stripping the footnote markers both become
#  "Influenza " and Excel naturally collapses them to one shared string.)

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    if ($used -eq $null) { continue }

    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $startRow = $used.Row
    $startCol = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($startRow + $r, $startCol + $c)
            $orig = $cell.Text

            if ([string]::IsNullOrEmpty($orig)) { continue }

            $new = $orig -replace '\[\d+(,\s*\d+)*\]', ''
            $new = $new -replace "`r`n", ' '
            $new = $new -replace "`n", ' '

            if ($new -ne $orig) {
                $cell.Value = $new
            }
        }
    }
}
